# The deck's single Design ("Integral") is backed by ppt/theme/theme2.xml,
# which supplies the 12-slot theme colour scheme used by every slide.
# The target edit swaps that scheme back to the stock "Office" palette
# (the same 12 colours that ship in ppt/theme/theme1.xml, which backs the
# notes master). Do it the way PowerPoint itself records it: walk the
# SlideMaster's theme colour scheme and set each RGBColor's .RGB.
#
# COM packs colours as 0xBBGGRR, so convert each target RRGGBB hex value
# before assigning.

function Set-ThemeRGB($themeColors, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Item($index).RGB = $bgr
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office theme palette: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
Set-ThemeRGB $colors 1  "000000"
Set-ThemeRGB $colors 2  "FFFFFF"
Set-ThemeRGB $colors 3  "44546A"
Set-ThemeRGB $colors 4  "E7E6E6"
Set-ThemeRGB $colors 5  "5B9BD5"
Set-ThemeRGB $colors 6  "ED7D31"
Set-ThemeRGB $colors 7  "A5A5A5"
Set-ThemeRGB $colors 8  "FFC000"
Set-ThemeRGB $colors 9  "4472C4"
Set-ThemeRGB $colors 10 "70AD47"
Set-ThemeRGB $colors 11 "0563C1"
Set-ThemeRGB $colors 12 "954F72"
